# Update need_to_buy.xlsx data rows (A2:F15) per latest R data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 45957
$ws.Cells.Item(2, 2).Value = 8634.17395320538
$ws.Cells.Item(2, 3).Value = 7903.95542890017
$ws.Cells.Item(2, 4).Value = 12020
$ws.Cells.Item(2, 5).Value = 4092.67064706103
$ws.Cells.Item(2, 6).Value = -0.973913501616456

# Row 3
$ws.Cells.Item(3, 1).Value = 45958
$ws.Cells.Item(3, 2).Value = 8327.82636748438
$ws.Cells.Item(3, 3).Value = 7779.68177796371
$ws.Cells.Item(3, 4).Value = 3620
$ws.Cells.Item(3, 5).Value = 4172.2893334047
$ws.Cells.Item(3, 6).Value = 347.165462973684

# Row 4
$ws.Cells.Item(4, 1).Value = 45959
$ws.Cells.Item(4, 2).Value = 7845.07976736576
$ws.Cells.Item(4, 3).Value = 7438.44813412283
$ws.Cells.Item(4, 4).Value = 3620
$ws.Cells.Item(4, 5).Value = 3853.3426395133
$ws.Cells.Item(4, 6).Value = 319.657948901505

# Row 5
$ws.Cells.Item(5, 1).Value = 45960
$ws.Cells.Item(5, 2).Value = 7640.54326920728
$ws.Cells.Item(5, 3).Value = 7018.46649202591
$ws.Cells.Item(5, 4).Value = 3620
$ws.Cells.Item(5, 5).Value = 3723.16100970739
$ws.Cells.Item(5, 6).Value = 296.734479238887

# Row 6
$ws.Cells.Item(6, 1).Value = 45961
$ws.Cells.Item(6, 2).Value = 7108.06757337707
$ws.Cells.Item(6, 3).Value = 5920.13069257766
$ws.Cells.Item(6, 4).Value = 3620
$ws.Cells.Item(6, 5).Value = 3354.55946611564
$ws.Cells.Item(6, 6).Value = 235.612089945554

# Row 7
$ws.Cells.Item(7, 1).Value = 45962
$ws.Cells.Item(7, 2).Value = 2557.16655997904
$ws.Cells.Item(7, 3).Value = 3694.79277565356
$ws.Cells.Item(7, 4).Value = 11404
$ws.Cells.Item(7, 5).Value = 3797.23033989697
$ws.Cells.Item(7, 6).Value = -162.999036852061

# Row 8
$ws.Cells.Item(8, 1).Value = 45963
$ws.Cells.Item(8, 2).Value = 2557.16655997904
$ws.Cells.Item(8, 3).Value = 3681.50298571129
$ws.Cells.Item(8, 4).Value = 11404
$ws.Cells.Item(8, 5).Value = 3797.23033989697
$ws.Cells.Item(8, 6).Value = -163.552778099656

# Row 9
$ws.Cells.Item(9, 1).Value = 45964
$ws.Cells.Item(9, 2).Value = 8571.71877057272
$ws.Cells.Item(9, 3).Value = 6975.47321521983
$ws.Cells.Item(9, 4).Value = 11404
$ws.Cells.Item(9, 5).Value = 4686.00445217985
$ws.Cells.Item(9, 6).Value = 10.7282361416534

# Row 10
$ws.Cells.Item(10, 1).Value = 45965
$ws.Cells.Item(10, 2).Value = 8571.71877057272
$ws.Cells.Item(10, 3).Value = 7729.24765355756
$ws.Cells.Item(10, 4).Value = 11404
$ws.Cells.Item(10, 5).Value = 4686.00445217985
$ws.Cells.Item(10, 6).Value = 42.1355044057254

# Row 11
$ws.Cells.Item(11, 1).Value = 45966
$ws.Cells.Item(11, 2).Value = 8571.71877057272
$ws.Cells.Item(11, 3).Value = 7842.63908281105
$ws.Cells.Item(11, 4).Value = 11404
$ws.Cells.Item(11, 5).Value = 4686.00445217985
$ws.Cells.Item(11, 6).Value = 46.8601472912877

# Row 12
$ws.Cells.Item(12, 1).Value = 45967
$ws.Cells.Item(12, 2).Value = 8580.99785695336
$ws.Cells.Item(12, 3).Value = 7722.66078489514
$ws.Cells.Item(12, 4).Value = 11404
$ws.Cells.Item(12, 5).Value = 4688.09354609001
$ws.Cells.Item(12, 6).Value = 41.9480971243813

# Row 13
$ws.Cells.Item(13, 1).Value = 45968
$ws.Cells.Item(13, 2).Value = 8585.92244945998
$ws.Cells.Item(13, 3).Value = 6912.03373629372
$ws.Cells.Item(13, 4).Value = 11404
$ws.Cells.Item(13, 5).Value = 4719.61696895587
$ws.Cells.Item(13, 6).Value = 9.48544605206644

# Row 14
$ws.Cells.Item(14, 1).Value = 45969
$ws.Cells.Item(14, 2).Value = 3147.97579797074
$ws.Cells.Item(14, 3).Value = 4636.65858310418
$ws.Cells.Item(14, 4).Value = 11404
$ws.Cells.Item(14, 5).Value = 4404.39210708016
$ws.Cells.Item(14, 6).Value = -98.4562212423193

# Row 15
$ws.Cells.Item(15, 1).Value = 45970
$ws.Cells.Item(15, 2).Value = 3012.17181266658
$ws.Cells.Item(15, 3).Value = 4473.36716439019
$ws.Cells.Item(15, 4).Value = 11404
$ws.Cells.Item(15, 5).Value = 4396.2524692723
$ws.Cells.Item(15, 6).Value = -105.59918193073

Write-Host "Updated rows 2-15"
